$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8633707761764526
$ws.Range("B1").Value = 1.540903568267822
$ws.Range("C1").Value = 4.362637042999268
$ws.Range("D1").Value = 2.080730438232422
$ws.Range("E1").Value = 1.515071749687195
